# feat: add PIE and USPS_MNIST dataset
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- USOS_MNIST block (rows 59-67) ---
# String literals must be written in this exact order so the shared-string
# table gets the same index assignment as the target workbook.
$ws.Range("B59").Value = "USOS_MNIST"
$ws.Range("C59").Value = "gamma=0.01"
$ws.Range("D59").Value = "gamma=0.1"

$ws.Range("B60").Value = 0.89
$ws.Range("C60").Value = 0.88722222222222202

$ws.Range("B61").Value = 0.748
$ws.Range("C61").Value = 0.73550000000000004

$ws.Range("B67").Value = "average"
$ws.Range("B64").Value = "relative performance"

$ws.Range("C64").Formula = '=C60-$B60'
$ws.Range("C65").Formula = '=C61-$B61'
$ws.Range("C67").Formula = '=AVERAGE(C64:C65)'

# --- PIE block (rows 70-90) ---
$ws.Range("A70").Value = "PIE"
$ws.Range("B70").Value = "MEDA"

$ws.Range("A71").Value = "05->07"
$ws.Range("B71").Value = 0.39472068753836698

$ws.Range("A72").Value = "05->09"
$ws.Range("B72").Value = 0.44607843137254899

$ws.Range("A73").Value = "05->27"
$ws.Range("B73").Value = 0.64854310603784904

$ws.Range("A74").Value = "05->29"
$ws.Range("B74").Value = 0.340073529411765

$ws.Range("A75").Value = "07->05"
$ws.Range("B75").Value = 0.46278511404561801

$ws.Range("A76").Value = "07->09"
$ws.Range("B76").Value = 0.50796568627451

$ws.Range("A77").Value = "07->27"
$ws.Range("B77").Value = 0.71312706518474001

$ws.Range("A78").Value = "07->29"
$ws.Range("B78").Value = 0.375

$ws.Range("A79").Value = "09->05"
$ws.Range("B79").Value = 0.46878751500600202

$ws.Range("A80").Value = "09->07"
$ws.Range("B80").Value = 0.488029465930018

$ws.Range("A81").Value = "09->27"
$ws.Range("B81").Value = 0.72874737158305802

$ws.Range("A82").Value = "09->29"
$ws.Range("B82").Value = 0.45465686274509798

$ws.Range("A83").Value = "27->05"
$ws.Range("B83").Value = 0.70048019207683099

$ws.Range("A84").Value = "27->07"
$ws.Range("B84").Value = 0.74708410067526099

$ws.Range("A85").Value = "27->09"
$ws.Range("B85").Value = 0.82046568627451

$ws.Range("A86").Value = "27->29"
$ws.Range("B86").Value = 0.54901960784313697

$ws.Range("A87").Value = "29->05"
$ws.Range("B87").Value = 0.38955582232893199

$ws.Range("A88").Value = "29->07"
$ws.Range("B88").Value = 0.37323511356660499

$ws.Range("A89").Value = "29->09"
$ws.Range("B89").Value = 0.44546568627451

$ws.Range("A90").Value = "29->27"
$ws.Range("B90").Value = 0.53259237008110505

# --- Apply the distinctive font (size 12, black) used for most PIE labels ---
$fontCells = @("A74","A75","A76","A77","A79","A80","A81","A82","A83","A84","A85","A86","A87","A88","A89","A90")
foreach ($addr in $fontCells) {
    $rng = $ws.Range($addr)
    $rng.Font.Size = 12
    $rng.Font.Color = 0
}

# --- View state: select B90 and scroll so row 59 is at the top ---
$ws.Range("B90").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 59
$win.ScrollColumn = 1
